# "new list and end of level"
#
# Remove the following bullet paragraphs from the proposed-features list:
#   - More static objects (platforms, obstacles, etc.)
#   - More enemies
#   - End of level (finish point, completion message, maybe new level?)
#   - Textures (enemies, objects, ground, etc.)
#   - Health system (if you have a power-up and get hit, you lose it. If you
#     get hit without a power-up, you die)
#   - Lives (find extra lives)
#   - Game over (need to print a screen for x seconds or something, then reset)
#   - Make level longer in general
#   - Maybe high score list?
#
# while keeping:
#   - (title) LIST OF THINGS TO ACCOMPLISH FOR CS430 PROJECT:
#   - Power-ups (get bigger, enemies move slower, etc.)
#   - More stuff for Points system (item collection, finish level, lives
#     bonus when you finish, etc.)
#   - Sound effects
#   - (trailing blank paragraphs)

$d = $word.ActiveDocument

# Work from the bottom of the list upward so earlier deletions don't shift
# the paragraph indices of the blocks still to be removed.

# Block 3: "Lives (find extra lives)" .. "Maybe high score list?"
$startRange = $d.Paragraphs.Item(9).Range.Start
$endRange = $d.Paragraphs.Item(12).Range.End
$d.Range($startRange, $endRange).Delete()

# Block 2: "Health system (...)"
$startRange = $d.Paragraphs.Item(7).Range.Start
$endRange = $d.Paragraphs.Item(7).Range.End
$d.Range($startRange, $endRange).Delete()

# Block 1: "More static objects (...)" .. "Textures (...)"
$startRange = $d.Paragraphs.Item(2).Range.Start
$endRange = $d.Paragraphs.Item(5).Range.End
$d.Range($startRange, $endRange).Delete()
